$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047636403393069
$ws.Range("D2").Value = 1.045330786168474
$ws.Range("E2").Value = 1.054809629692547
$ws.Range("F2").Value = 1.065058965107623
$ws.Range("I2").Value = 1.037722866907686
$ws.Range("J2").Value = 1.052684416159849
$ws.Range("K2").Value = 1.048099327578221
$ws.Range("L2").Value = 1.057551816867986
$ws.Range("M2").Value = 1.067773227213547
$ws.Range("N2").Value = 1.054179348775386

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049099351950871
$ws.Range("D3").Value = 1.046387318368268
$ws.Range("E3").Value = 1.056122647936853
$ws.Range("F3").Value = 1.066484358712065
$ws.Range("I3").Value = 1.038030198907875
$ws.Range("J3").Value = 1.05379339559007
$ws.Range("K3").Value = 1.04896664022097
$ws.Range("L3").Value = 1.058676876088282
$ws.Range("M3").Value = 1.069012420843264
$ws.Range("N3").Value = 1.0552899030836

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050044847489913
$ws.Range("D4").Value = 1.047069782557128
$ws.Range("E4").Value = 1.056971507606922
$ws.Range("F4").Value = 1.067406054740654
$ws.Range("I4").Value = 1.038227208199964
$ws.Range("J4").Value = 1.054509440044578
$ws.Range("K4").Value = 1.049526071514289
$ws.Range("L4").Value = 1.059403569294328
$ws.Range("M4").Value = 1.069813107880382
$ws.Range("N4").Value = 1.056006964403364

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050442069588473
$ws.Range("D5").Value = 1.047356410861164
$ws.Range("E5").Value = 1.057328193217293
$ws.Range("F5").Value = 1.067793390398083
$ws.Range("I5").Value = 1.03830958818263
$ws.Range("J5").Value = 1.054810101067374
$ws.Range("K5").Value = 1.049760834140149
$ws.Range("L5").Value = 1.059708765382906
$ws.Range("M5").Value = 1.070149445129731
$ws.Range("N5").Value = 1.056308052399299

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050508749614862
$ws.Range("D6").Value = 1.047404520718297
$ws.Range("E6").Value = 1.057388072177446
$ws.Range("F6").Value = 1.067858417386942
$ws.Range("I6").Value = 1.038323394222663
$ws.Range("J6").Value = 1.054860562156471
$ws.Range("K6").Value = 1.049800227152097
$ws.Range("L6").Value = 1.059759991347268
$ws.Range("M6").Value = 1.070205901844755
$ws.Range("N6").Value = 1.056358585148931

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050050156221459
$ws.Range("D7").Value = 1.047073613593005
$ws.Range("E7").Value = 1.056976274339197
$ws.Range("F7").Value = 1.067411230903337
$ws.Range("I7").Value = 1.038228310703224
$ws.Range("J7").Value = 1.054513458916888
$ws.Range("K7").Value = 1.049529210079149
$ws.Range("L7").Value = 1.059407648535325
$ws.Range("M7").Value = 1.069817603096465
$ws.Range("N7").Value = 1.056010988982934

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048131051481598
$ws.Range("D8").Value = 1.045688093203618
$ws.Range("E8").Value = 1.055253528725512
$ws.Range("F8").Value = 1.065540816925616
$ws.Range("I8").Value = 1.03782711613377
$ws.Range("J8").Value = 1.053059522387852
$ws.Range("K8").Value = 1.048392810119511
$ws.Range("L8").Value = 1.057932306745739
$ws.Range("M8").Value = 1.068192260311543
$ws.Range("N8").Value = 1.054554987697258

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044740406404847
$ws.Range("D9").Value = 1.043237411787186
$ws.Range("E9").Value = 1.052211864693109
$ws.Range("F9").Value = 1.06223985702887
$ws.Range("I9").Value = 1.037105890117173
$ws.Range("J9").Value = 1.050485509641142
$ws.Range("K9").Value = 1.046376559551044
$ws.Range("L9").Value = 1.055322448667743
$ws.Range("M9").Value = 1.065319147685586
$ws.Range("N9").Value = 1.051977319557204

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042473570048405
$ws.Range("D10").Value = 1.041597202467081
$ws.Range("E10").Value = 1.050179774671977
$ws.Range("F10").Value = 1.060035490652371
$ws.Range("I10").Value = 1.036615390011907
$ws.Range("J10").Value = 1.048761156171883
$ws.Range("K10").Value = 1.045022921375665
$ws.Range("L10").Value = 1.05357546400264
$ws.Range("M10").Value = 1.063397355435017
$ws.Range("N10").Value = 1.050250517308234

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041490396419312
$ws.Range("D11").Value = 1.040885399877416
$ws.Range("E11").Value = 1.049298766012557
$ws.Range("F11").Value = 1.059080017692314
$ws.Range("I11").Value = 1.036400681339825
$ws.Range("J11").Value = 1.048012449150002
$ws.Range("K11").Value = 1.044434486428473
$ws.Range("L11").Value = 1.052817258989448
$ws.Range("M11").Value = 1.062563616944062
$ws.Range("N11").Value = 1.049500747036503

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041124951047703
$ws.Range("D12").Value = 1.040620762691742
$ws.Range("E12").Value = 1.048971348921355
$ws.Range("F12").Value = 1.0587249606255
$ws.Range("I12").Value = 1.036320578798025
$ws.Range("J12").Value = 1.047734032378287
$ws.Range("K12").Value = 1.044215565589775
$ws.Range("L12").Value = 1.052535359653276
$ws.Range("M12").Value = 1.062253684636239
$ws.Range("N12").Value = 1.049221934881037

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041203351799538
$ws.Range("D13").Value = 1.040677539289953
$ws.Range("E13").Value = 1.0490415888441
$ws.Range("F13").Value = 1.058801128548262
$ws.Range("I13").Value = 1.036337776940544
$ws.Range("J13").Value = 1.047793768003985
$ws.Range("K13").Value = 1.044262540755037
$ws.Range("L13").Value = 1.052595840243919
$ws.Range("M13").Value = 1.06232017733662
$ws.Range("N13").Value = 1.049281755338176

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04146019375436
$ws.Range("D14").Value = 1.040863529847991
$ws.Range("E14").Value = 1.049271705141608
$ws.Range("F14").Value = 1.059050671700036
$ws.Range("I14").Value = 1.03639406719275
$ws.Range("J14").Value = 1.047989441563473
$ws.Range("K14").Value = 1.044416397526177
$ws.Range("L14").Value = 1.052793962605954
$ws.Range("M14").Value = 1.062538002868901
$ws.Range("N14").Value = 1.049477706776563

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041618409079028
$ws.Range("D15").Value = 1.040978092566324
$ws.Range("E15").Value = 1.049413464516693
$ws.Range("F15").Value = 1.059204403153642
$ws.Range("I15").Value = 1.036428703018637
$ws.Range("J15").Value = 1.048109960772874
$ws.Range("K15").Value = 1.044511147268128
$ws.Range("L15").Value = 1.052915996624636
$ws.Range("M15").Value = 1.062672179765852
$ws.Range("N15").Value = 1.049598397137065

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042538785283772
$ws.Range("D16").Value = 1.041644408767887
$ws.Range("E16").Value = 1.050238220627252
$ws.Range("F16").Value = 1.06009888130485
$ws.Range("I16").Value = 1.036629590510193
$ws.Range("J16").Value = 1.048810801713038
$ws.Range("K16").Value = 1.045061925030565
$ws.Range("L16").Value = 1.053625746256004
$ws.Range("M16").Value = 1.063452653918539
$ws.Range("N16").Value = 1.050300233351751

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043115674183209
$ws.Range("D17").Value = 1.042061944863153
$ws.Range("E17").Value = 1.050755269587606
$ws.Range("F17").Value = 1.060659700162744
$ws.Range("I17").Value = 1.036754979831339
$ws.Range("J17").Value = 1.049249867797479
$ws.Range("K17").Value = 1.045406794422234
$ws.Range("L17").Value = 1.054070481322828
$ws.Range("M17").Value = 1.063941794980664
$ws.Range("N17").Value = 1.050739922960394

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043452008079985
$ws.Range("D18").Value = 1.04230533435141
$ws.Range("E18").Value = 1.051056749742758
$ws.Range("F18").Value = 1.060986723378371
$ws.Range("I18").Value = 1.03682789368402
$ws.Range("J18").Value = 1.049505770142796
$ws.Range("K18").Value = 1.045607729210214
$ws.Range("L18").Value = 1.054329719368546
$ws.Range("M18").Value = 1.064226949696442
$ws.Range("N18").Value = 1.050996188716393

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043566663081436
$ws.Range("D19").Value = 1.042388298186754
$ws.Range("E19").Value = 1.051159528925611
$ws.Range("F19").Value = 1.061098214275046
$ws.Range("I19").Value = 1.036852717548117
$ws.Range("J19").Value = 1.049592992902132
$ws.Range("K19").Value = 1.045676205304962
$ws.Range("L19").Value = 1.054418084444167
$ws.Range("M19").Value = 1.0643241543122
$ws.Range("N19").Value = 1.051083535342052

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043053795618494
$ws.Range("D20").Value = 1.042017162972228
$ws.Range("E20").Value = 1.050699806109441
$ws.Range("F20").Value = 1.060599539266612
$ws.Range("I20").Value = 1.036741549886643
$ws.Range("J20").Value = 1.04920278062793
$ws.Range("K20").Value = 1.045369816161966
$ws.Range("L20").Value = 1.054022782916484
$ws.Range("M20").Value = 1.063889330640972
$ws.Range("N20").Value = 1.050692768921663

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041384567189075
$ws.Range("D21").Value = 1.040808766981053
$ws.Range("E21").Value = 1.049203946402731
$ws.Range("F21").Value = 1.058977191705993
$ws.Range("I21").Value = 1.036377500797664
$ws.Range("J21").Value = 1.047931829280374
$ws.Range("K21").Value = 1.044371100235059
$ws.Range("L21").Value = 1.052735627956961
$ws.Range("M21").Value = 1.062473865481311
$ws.Range("N21").Value = 1.049420012677413

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040333601463464
$ws.Range("D22").Value = 1.040047598023059
$ws.Range("E22").Value = 1.0482624462704
$ws.Range("F22").Value = 1.057956275511182
$ws.Range("I22").Value = 1.036146581747522
$ws.Range("J22").Value = 1.04713091426854
$ws.Range("K22").Value = 1.043741141356388
$ws.Range("L22").Value = 1.05192478789389
$ws.Range("M22").Value = 1.061582487603662
$ws.Range("N22").Value = 1.048617960274396

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040890878847859
$ws.Range("D23").Value = 1.040451242523284
$ws.Range("E23").Value = 1.048761649496611
$ws.Range("F23").Value = 1.058497568266505
$ws.Range("I23").Value = 1.036269189037935
$ws.Range("J23").Value = 1.047555668612971
$ws.Range("K23").Value = 1.044075287995325
$ws.Range("L23").Value = 1.052354778744575
$ws.Range("M23").Value = 1.062055160331605
$ws.Range("N23").Value = 1.049043317818715

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043081756359567
$ws.Range("D24").Value = 1.042037398451889
$ws.Range("E24").Value = 1.050724867994938
$ws.Range("F24").Value = 1.060626723673497
$ws.Range("I24").Value = 1.03674761899165
$ws.Range("J24").Value = 1.049224057903465
$ws.Range("K24").Value = 1.045386525731342
$ws.Range("L24").Value = 1.054044336293772
$ws.Range("M24").Value = 1.063913037490822
$ws.Range("N24").Value = 1.05071407641337

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045618069311216
$ws.Range("D25").Value = 1.043872086987068
$ws.Range("E25").Value = 1.052998945813616
$ws.Range("F25").Value = 1.063093866769847
$ws.Range("I25").Value = 1.037294044453012
$ws.Range("J25").Value = 1.051152403628782
$ws.Range("K25").Value = 1.046899462227065
$ws.Range("L25").Value = 1.055998386894671
$ws.Range("M25").Value = 1.066063019362498
$ws.Range("N25").Value = 1.052645160610801
